$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 176
$ws1.Range("F6").Value = 673
$ws1.Range("F8").Value = 483
$ws1.Range("F10").Value = 527
$ws1.Range("F11").Value = 405
$ws1.Range("F12").Value = 68
$ws1.Range("F14").Value = 117
$ws1.Range("F15").Value = 200

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 22

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 1824

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1824
$ws4.Range("F10").Value = 22
$ws4.Range("F12").Value = 176
$ws4.Range("F15").Value = 673
$ws4.Range("F19").Value = 483
$ws4.Range("F22").Value = 527
$ws4.Range("F24").Value = 405
$ws4.Range("F25").Value = 68
$ws4.Range("F29").Value = 117
$ws4.Range("F35").Value = 200
